$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 167.57143
$ws.Range("I6").Value = 190.5
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 571.5
$ws.Range("L6").Value = 90
$ws.Range("M6").Value = -459.5
$ws.Range("N6").Value = -314
# Row 64
$ws.Range("H64").Value = 6857.143
$ws.Range("I64").Value = 5750
$ws.Range("J64").Value = 8333.333000000001
$ws.Range("K64").Value = 5750
$ws.Range("L64").Value = 8333.333000000001
$ws.Range("M64").Value = -5502
$ws.Range("N64").Value = -8829.333000000001
# Row 67
$ws.Range("H67").Value = 6857.143
$ws.Range("I67").Value = 5750
$ws.Range("J67").Value = 8333.333000000001
$ws.Range("K67").Value = 5750
$ws.Range("L67").Value = 8333.333000000001
$ws.Range("M67").Value = -4892
$ws.Range("N67").Value = -10049.333
# Row 70
$ws.Range("H70").Value = 3590.818
$ws.Range("I70").Value = 3388.889
$ws.Range("K70").Value = 10166.667
$ws.Range("M70").Value = -9896.667000000001
# Row 73
$ws.Range("H73").Value = 3590.818
$ws.Range("I73").Value = 3388.889
$ws.Range("K73").Value = 10166.667
$ws.Range("M73").Value = -9230.667000000001
# Row 88
$ws.Range("H88").Value = 1489.5238
$ws.Range("J88").Value = 1546.421
$ws.Range("L88").Value = 1546.421
$ws.Range("N88").Value = -2358.421
# Row 91
$ws.Range("H91").Value = 1489.5238
$ws.Range("J91").Value = 1546.421
$ws.Range("L91").Value = 1546.421
$ws.Range("N91").Value = -4354.421
# Row 100
$ws.Range("H100").Value = 4001.8
$ws.Range("I100").Value = 3336.6667
$ws.Range("J100").Value = 4999.5
$ws.Range("K100").Value = 3336.6667
$ws.Range("L100").Value = 4999.5
$ws.Range("M100").Value = -2795.6667
$ws.Range("N100").Value = -6081.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3482
$ws.Range("I2").Value = 3482
$ws.Range("K2").Value = 3482
$ws.Range("M2").Value = -3369
# Row 4
$ws.Range("H4").Value = 901
$ws.Range("I4").Value = 901
$ws.Range("K4").Value = 901
$ws.Range("M4").Value = -785
# Row 29
$ws.Range("H29").Value = 355.25
$ws.Range("I29").Value = 355.25
$ws.Range("K29").Value = 355.25
$ws.Range("M29").Value = -47.25
# Row 74
$ws.Range("H74").Value = 1608.75
$ws.Range("I74").Value = 1541
$ws.Range("K74").Value = 1541
$ws.Range("M74").Value = -667
# Row 77
$ws.Range("H77").Value = 1608.75
$ws.Range("I77").Value = 1541
$ws.Range("K77").Value = 7705
$ws.Range("M77").Value = -3337
# Row 97
$ws.Range("H97").Value = 66670500
$ws.Range("I97").Value = 66670500
$ws.Range("K97").Value = 66670500
$ws.Range("M97").Value = -66670004
# Row 116
$ws.Range("H116").Value = 3482
$ws.Range("I116").Value = 3482
$ws.Range("K116").Value = 3482
$ws.Range("M116").Value = -1188
# Row 122
$ws.Range("H122").Value = 5510.0835
$ws.Range("I122").Value = 2437
$ws.Range("J122").Value = 8583.166999999999
$ws.Range("K122").Value = 7311
$ws.Range("L122").Value = 25749.501
$ws.Range("M122").Value = -4861
$ws.Range("N122").Value = -30649.501

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3482
$ws.Range("I3").Value = 3482
$ws.Range("K3").Value = 3482
$ws.Range("M3").Value = -3368
# Row 33
$ws.Range("H33").Value = 47499.75
$ws.Range("I33").Value = 39999.5
$ws.Range("K33").Value = 39999.5
$ws.Range("M33").Value = -39663.5
# Row 99
$ws.Range("H99").Value = 2929.6667
$ws.Range("I99").Value = 2929.6667
$ws.Range("K99").Value = 2929.6667
$ws.Range("M99").Value = -1431.6667
# Row 105
$ws.Range("H105").Value = 1633.3334
$ws.Range("I105").Value = 1633.3334
$ws.Range("K105").Value = 1633.3334
$ws.Range("M105").Value = 113.6666
# Row 107
$ws.Range("H107").Value = 38490.363
$ws.Range("I107").Value = 41839.7
$ws.Range("J107").Value = 4997
$ws.Range("K107").Value = 41839.7
$ws.Range("L107").Value = 4997
$ws.Range("M107").Value = -39919.7
$ws.Range("N107").Value = -8837

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 600
$ws.Range("I4").Value = 600
$ws.Range("K4").Value = 600
$ws.Range("M4").Value = -488
# Row 12
$ws.Range("H12").Value = 1667502.1
$ws.Range("I12").Value = 2000442.6
$ws.Range("K12").Value = 2000442.6
$ws.Range("M12").Value = -2000272.6
# Row 35
$ws.Range("H35").Value = 998.2
$ws.Range("I35").Value = 998.2
$ws.Range("K35").Value = 998.2
$ws.Range("M35").Value = -704.2
# Row 48
$ws.Range("H48").Value = 43666.668
$ws.Range("J48").Value = 43666.668
$ws.Range("L48").Value = 43666.668
$ws.Range("N48").Value = -44618.668
# Row 122
$ws.Range("H122").Value = 2630.8333
$ws.Range("I122").Value = 2630.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7892.499899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5442.499899999999
$ws.Range("N122").ClearContents()
# Row 141
$ws.Range("H141").Value = 813608.8
$ws.Range("J141").Value = 813608.8
$ws.Range("L141").Value = 813608.8
$ws.Range("N141").Value = -823968.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 197.5
$ws.Range("I23").Value = 163.33333
$ws.Range("K23").Value = 489.99999
$ws.Range("M23").Value = -254.99999
# Row 40
$ws.Range("H40").Value = 172.16667
$ws.Range("I40").Value = 46.6
$ws.Range("J40").Value = 800
$ws.Range("K40").Value = 186.4
$ws.Range("L40").Value = 3200
$ws.Range("M40").Value = -117.4
$ws.Range("N40").Value = -3338
# Row 55
$ws.Range("H55").Value = 4019.6885
$ws.Range("J55").Value = 4130.3896
$ws.Range("L55").Value = 12391.1688
$ws.Range("N55").Value = -12745.1688
# Row 98
$ws.Range("H98").Value = 1009.625
$ws.Range("J98").Value = 1416
$ws.Range("L98").Value = 4248
$ws.Range("N98").Value = -7244
# Row 137
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Range("H53").Value = 12833.333
$ws.Range("J53").Value = 14250
$ws.Range("L53").Value = 14250
$ws.Range("N53").Value = -15512
# Row 80
$ws.Range("H80").Value = 12049.833
$ws.Range("I80").Value = 16075
$ws.Range("J80").Value = 3999.5
$ws.Range("K80").Value = 16075
$ws.Range("L80").Value = 3999.5
$ws.Range("M80").Value = -15077
$ws.Range("N80").Value = -5995.5
# Row 83
$ws.Range("H83").Value = 12049.833
$ws.Range("I83").Value = 16075
$ws.Range("J83").Value = 3999.5
$ws.Range("K83").Value = 80375
$ws.Range("L83").Value = 19997.5
$ws.Range("M83").Value = -75383
$ws.Range("N83").Value = -29981.5
# Row 97
$ws.Range("H97").Value = 2807.5
$ws.Range("I97").Value = 2807.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2807.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2311.5
$ws.Range("N97").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 359.33334
$ws.Range("I9").Value = 39
$ws.Range("K9").Value = 39
$ws.Range("M9").Value = 185
# Row 30
$ws.Range("H30").Value = 950
$ws.Range("I30").Value = 950
$ws.Range("K30").Value = 950
$ws.Range("M30").Value = -842
# Row 100
$ws.Range("H100").Value = 3170.6
$ws.Range("I100").Value = 3170.6
$ws.Range("K100").Value = 3170.6
$ws.Range("M100").Value = -2629.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 33776.75
$ws.Range("I61").Value = 19525
$ws.Range("J61").Value = 48028.5
$ws.Range("K61").Value = 19525
$ws.Range("L61").Value = 48028.5
$ws.Range("M61").Value = -19233
$ws.Range("N61").Value = -48612.5
# Row 81
$ws.Range("H81").Value = 2404.8
$ws.Range("I81").Value = 2404.8
$ws.Range("K81").Value = 4809.6
$ws.Range("M81").Value = -3748.6
# Row 84
$ws.Range("H84").Value = 2404.8
$ws.Range("I84").Value = 2404.8
$ws.Range("K84").Value = 24048
$ws.Range("M84").Value = -18744
# Row 107
$ws.Range("H107").Value = 550
$ws.Range("I107").Value = 433.33334
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1300.00002
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 619.9999800000001
$ws.Range("N107").Value = -6540
